# Updated cryptos list with refreshed prices / 1h volume percentages.
# Cells in column D hold price strings that look like plain numbers
# (e.g. "0.9995", "1.000", "0.05330"); Excel would silently coerce these
# to numeric values (losing trailing zeros / exact formatting) unless the
# cell is explicitly forced to Text format first. Values that use the
# site's "thousands" dotted notation (e.g. "26.488.91") are left alone
# since Excel already keeps them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.488.91"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "1.727.23"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.48"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4797"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2684"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06225"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "1.727.27"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07144"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.70"
$ws.Range("E12").Value = "  +3.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6192"
$ws.Range("E13").Value = "  +5.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.516"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.14"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "26.496.14"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006939"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.67"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "1.949.89"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.536"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.961"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.295"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.41"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.34"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.803"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.76"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.977"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08020"
$ws.Range("E31").Value = "  +3.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.722"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04568"
$ws.Range("E33").Value = "  +3.59%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9994"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.618"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6366"
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9890"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9363"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.113"
$ws.Range("E39").Value = "  +11.38%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.412"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.96"
$ws.Range("E41").Value = "  -7.09%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01502"
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.697"
$ws.Range("E44").Value = "  +9.26%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3898"
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.939"
$ws.Range("E46").Value = "  +11.37%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1189"
$ws.Range("E47").Value = "  +3.85%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05330"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.04"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.873"
$ws.Range("E50").Value = "  +2.62%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.266"
$ws.Range("E51").Value = "  +3.93%  "
